$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update "last updated" timestamp (A1)
$ws.Range("A1").Value = "Datos actualizados a 20 de Mayo de 2020 a las 21:05"

# Re-rank rows whose country order changed due to new totals (keeps same shared-string slot index)
$ws.Range("A37").Value = "Sudafrica"
$ws.Range("A38").Value = "Kuwait"
$ws.Range("A39").Value = "Rumania"
$ws.Range("A110").Value = "Mali"
$ws.Range("A111").Value = "Republica de Chipre"
$ws.Range("A112").Value = "Niger"
$ws.Range("A196").Value = "Santa Lucia"
$ws.Range("A197").Value = "Belice"
$ws.Range("A209").Value = "Montserrat"
$ws.Range("A210").Value = "Groenlandia"

# Update numeric stats per diff
$ws.Range("B4").Value = 1579387
$ws.Range("C4").Value = 8804
$ws.Range("E4").Value = 1120537
$ws.Range("G4").Value = 634
$ws.Range("H4").Value = 94167
$ws.Range("B10").Value = 181575
$ws.Range("C10").Value = 766
$ws.Range("D10").Value = 63354
$ws.Range("E10").Value = 90089
$ws.Range("G10").Value = 110
$ws.Range("H10").Value = 28132
$ws.Range("B11").Value = 178344
$ws.Range("C11").Value = 517
$ws.Range("E11").Value = 13214
$ws.Range("G11").Value = 37
$ws.Range("H11").Value = 8230
$ws.Range("B14").Value = 112012
$ws.Range("C14").Value = 5537
$ws.Range("E14").Value = 63156
$ws.Range("D33").Value = 21060
$ws.Range("E33").Value = 1684
$ws.Range("B37").Value = 18003
$ws.Range("C37").Value = 803
$ws.Range("D37").Value = 8950
$ws.Range("E37").Value = 8714
$ws.Range("G37").Value = 27
$ws.Range("H37").Value = 339
$ws.Range("B38").Value = 17568
$ws.Range("C38").Value = 804
$ws.Range("D38").Value = 4885
$ws.Range("E38").Value = 12559
$ws.Range("G38").Value = 3
$ws.Range("H38").Value = 124
$ws.Range("B39").Value = 17387
$ws.Range("C39").Value = 196
$ws.Range("D39").Value = 10356
$ws.Range("E39").Value = 5884
$ws.Range("G39").Value = 10
$ws.Range("H39").Value = 1147
$ws.Range("B55").Value = 7888
$ws.Range("C55").Value = 356
$ws.Range("D55").Value = 3568
$ws.Range("E55").Value = 4308
$ws.Range("B110").Value = 931
$ws.Range("C110").Value = 30
$ws.Range("D110").Value = 543
$ws.Range("E110").Value = 333
$ws.Range("G110").Value = 2
$ws.Range("H110").Value = 55
$ws.Range("B111").Value = 922
$ws.Range("C111").Value = 4
$ws.Range("D111").Value = 516
$ws.Range("E111").Value = 389
$ws.Range("H111").Value = 17
$ws.Range("B112").Value = 914
$ws.Range("D112").Value = 734
$ws.Range("E112").Value = 125
$ws.Range("H112").Value = 55
$ws.Range("B113").Value = 897
$ws.Range("C113").Value = 15
$ws.Range("D113").Value = 582
$ws.Range("E113").Value = 305
$ws.Range("D196").Value = 18
$ws.Range("H196").Value = 0
$ws.Range("D197").Value = 16
$ws.Range("H197").Value = 2
$ws.Range("D209").Value = 10
$ws.Range("H209").Value = 1
$ws.Range("D210").Value = 11
$ws.Range("H210").Value = 0
